$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Sheets.Item("ALC")
$ws.Range("H18").Value = 1825.5
$ws.Range("I18").Value = 1825.5
$ws.Range("K18").Value = 1825.5
$ws.Range("M18").Value = -1541.5
$ws.Range("H19").Value = 899.875
$ws.Range("J19").Value = 1734
$ws.Range("L19").Value = 1734
$ws.Range("N19").Value = -2084
$ws.Range("H88").Value = 3000
$ws.Range("J88").Value = 3000
$ws.Range("L88").Value = 3000
$ws.Range("N88").Value = -3812
$ws.Range("H91").Value = 3000
$ws.Range("J91").Value = 3000
$ws.Range("L91").Value = 3000
$ws.Range("N91").Value = -5808
$ws.Range("H100").Value = 2005286.2
$ws.Range("I100").Value = 2505607.8
$ws.Range("J100").Value = 4000
$ws.Range("K100").Value = 2505607.8
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -2505066.8
$ws.Range("N100").Value = -5082
$ws.Range("H107").Value = 1952.5
$ws.Range("I107").Value = 1952.5
$ws.Range("K107").Value = 1952.5
$ws.Range("M107").Value = -32.5
$ws.Range("H138").Value = 8223.538
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 8223.538
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 24670.614
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -34950.614

# --- Sheet: ARM ---
$ws = $wb.Sheets.Item("ARM")
$ws.Range("H2").Value = 974.8570999999999
$ws.Range("I2").Value = 961.2
$ws.Range("K2").Value = 961.2
$ws.Range("M2").Value = -848.2
$ws.Range("H45").Value = 1932.2667
$ws.Range("I45").Value = 1932.2667
$ws.Range("K45").Value = 1932.2667
$ws.Range("M45").Value = -1555.2667
$ws.Range("H61").Value = 6921
$ws.Range("I61").Value = 6921
$ws.Range("K61").Value = 6921
$ws.Range("M61").Value = -6709
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H116").Value = 974.8570999999999
$ws.Range("I116").Value = 961.2
$ws.Range("K116").Value = 961.2
$ws.Range("M116").Value = 1332.8
$ws.Range("H122").Value = 24998.572
$ws.Range("I122").Value = 24998.334
$ws.Range("K122").Value = 74995.00199999999
$ws.Range("M122").Value = -72545.00199999999
$ws.Range("H132").Value = 3916
$ws.Range("I132").Value = 2833
$ws.Range("K132").Value = 8499
$ws.Range("M132").Value = -5969
$ws.Range("H136").Value = 6921
$ws.Range("I136").Value = 6921
$ws.Range("K136").Value = 20763
$ws.Range("M136").Value = -18213

# --- Sheet: BSM ---
$ws = $wb.Sheets.Item("BSM")
$ws.Range("H3").Value = 974.8570999999999
$ws.Range("I3").Value = 961.2
$ws.Range("K3").Value = 961.2
$ws.Range("M3").Value = -847.2
$ws.Range("H80").Value = 639.4545000000001
$ws.Range("J80").Value = 390
$ws.Range("L80").Value = 390
$ws.Range("N80").Value = -2386
$ws.Range("H81").Value = 73127.86
$ws.Range("J81").Value = 73127.86
$ws.Range("L81").Value = 73127.86
$ws.Range("N81").Value = -75249.86
$ws.Range("H83").Value = 639.4545000000001
$ws.Range("J83").Value = 390
$ws.Range("L83").Value = 1950
$ws.Range("N83").Value = -11934
$ws.Range("H84").Value = 73127.86
$ws.Range("J84").Value = 73127.86
$ws.Range("L84").Value = 219383.58
$ws.Range("N84").Value = -229991.58
$ws.Range("H130").Value = 43780
$ws.Range("J130").Value = 43780
$ws.Range("L130").Value = 43780
$ws.Range("N130").Value = -53820

# --- Sheet: CRP ---
$ws = $wb.Sheets.Item("CRP")
$ws.Range("H31").Value = 7249.5
$ws.Range("I31").Value = 2000
$ws.Range("J31").Value = 12499
$ws.Range("K31").Value = 2000
$ws.Range("L31").Value = 12499
$ws.Range("M31").Value = -1705
$ws.Range("N31").Value = -13089
$ws.Range("H34").Value = 7249.5
$ws.Range("I34").Value = 2000
$ws.Range("J34").Value = 12499
$ws.Range("K34").Value = 2000
$ws.Range("L34").Value = 12499
$ws.Range("M34").Value = -1798
$ws.Range("N34").Value = -12903
$ws.Range("H62").Value = 7199.6665
$ws.Range("J62").Value = 7199
$ws.Range("L62").Value = 7199
$ws.Range("N62").Value = -8447
$ws.Range("H65").Value = 7199.6665
$ws.Range("J65").Value = 7199
$ws.Range("L65").Value = 35995
$ws.Range("N65").Value = -42235
$ws.Range("H86").Value = 7481.125
$ws.Range("I86").Value = 8000.5
$ws.Range("J86").Value = 6961.75
$ws.Range("K86").Value = 8000.5
$ws.Range("L86").Value = 6961.75
$ws.Range("M86").Value = -6877.5
$ws.Range("N86").Value = -9207.75
$ws.Range("H89").Value = 7481.125
$ws.Range("I89").Value = 8000.5
$ws.Range("J89").Value = 6961.75
$ws.Range("K89").Value = 40002.5
$ws.Range("L89").Value = 34808.75
$ws.Range("M89").Value = -34386.5
$ws.Range("N89").Value = -46040.75
$ws.Range("H105").Value = 3982.3333
$ws.Range("I105").Value = 5465.6665
$ws.Range("J105").Value = 2499
$ws.Range("K105").Value = 5465.6665
$ws.Range("L105").Value = 2499
$ws.Range("M105").Value = -3718.6665
$ws.Range("N105").Value = -5993
$ws.Range("H107").Value = 2051.4
$ws.Range("I107").Value = 1419.3334
$ws.Range("J107").Value = 2999.5
$ws.Range("K107").Value = 1419.3334
$ws.Range("L107").Value = 2999.5
$ws.Range("M107").Value = 500.6666
$ws.Range("N107").Value = -6839.5

# --- Sheet: GSM ---
$ws = $wb.Sheets.Item("GSM")
$ws.Range("H126").Value = 5832.6665
$ws.Range("I126").Value = 3749.5
$ws.Range("K126").Value = 11248.5
$ws.Range("M126").Value = -8778.5

# --- Sheet: LTW ---
$ws = $wb.Sheets.Item("LTW")
$ws.Range("H61").Value = 3686.75
$ws.Range("I61").Value = 3665.6667
$ws.Range("K61").Value = 3665.6667
$ws.Range("M61").Value = -3463.6667
$ws.Range("H113").Value = 3686.75
$ws.Range("I113").Value = 3665.6667
$ws.Range("K113").Value = 3665.6667
$ws.Range("M113").Value = -1495.6667
$ws.Range("H136").Value = 2199.9
$ws.Range("I136").Value = 2250
$ws.Range("K136").Value = 6750
$ws.Range("M136").Value = -4200

# --- Sheet: WVR ---
$ws = $wb.Sheets.Item("WVR")
$ws.Range("H107").Value = 1114
$ws.Range("I107").Value = 1072.5
$ws.Range("K107").Value = 3217.5
$ws.Range("M107").Value = -1297.5
